# Auto-generated edit script applying the cryptos.xlsx update
# (crypto price/volume refresh + two coin-row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.297.75"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "3.070.49"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'210.72"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'626.50"
$ws.Range("E6").Value = "  +1.80%  "
$ws.Range("D7").Value = "'0.370"
$ws.Range("E7").Value = "  -5.40%  "
$ws.Range("D8").Value = "'0.791"
$ws.Range("E8").Value = "  +14.54%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "3.070.11"
$ws.Range("E10").Value = "  -2.94%  "
$ws.Range("D11").Value = "'0.580"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "'0.178"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -5.67%  "
$ws.Range("D14").Value = "'5.27"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "87.471.74"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "3.656.19"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "'31.42"
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "3.095.60"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "'3.38"
$ws.Range("E19").Value = "  +3.27%  "
$ws.Range("D20").Value = "'0.0000210"
$ws.Range("E20").Value = "  +7.14%  "
$ws.Range("D21").Value = "'13.08"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").Value = "'416.85"
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("D23").Value = "'8.25"
$ws.Range("E23").Value = "  -3.84%  "
$ws.Range("D24").Value = "'4.80"
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("D25").Value = "'5.39"
$ws.Range("E25").Value = "  +5.21%  "
$ws.Range("D26").Value = "'82.88"
$ws.Range("E26").Value = "  +10.56%  "
$ws.Range("D27").Value = "'11.21"
$ws.Range("E27").Value = "  -3.33%  "
$ws.Range("D28").Value = "3.258.33"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'0.150"
$ws.Range("E31").Value = "  -9.63%  "
$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("D33").Value = "'3.69"
$ws.Range("E33").Value = "  -8.71%  "
$ws.Range("D34").Value = "'495.38"
$ws.Range("E34").Value = "  -6.97%  "
$ws.Range("D35").Value = "'0.143"
$ws.Range("E35").Value = "  +12.97%  "
$ws.Range("D36").Value = "'6.69"
$ws.Range("E36").Value = "  -3.97%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.24"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'22.07"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").Value = "'22.15"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "'0.360"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("D44").Value = "'1.81"
$ws.Range("E44").Value = "  -5.92%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'146.61"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.134"
$ws.Range("E46").Value = "  +8.34%  "
$ws.Range("D47").Value = "'43.55"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "'0.0648"
$ws.Range("E48").Value = "  +9.83%  "
$ws.Range("D49").Value = "'158.90"
$ws.Range("E49").Value = "  -7.56%  "
$ws.Range("D50").Value = "'0.707"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "'1.17"
$ws.Range("E51").Value = "  -4.64%  "
